$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Updated fitting parameters
$ws.Range("J2").Value = 0.0337
$ws.Range("K2").Value = 0.23866

# Move selection to K3 (ready to run detail tests)
$ws.Activate()
$ws.Range("K3").Select()
